$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 39 (pushes old rows 39..169 down to 40..170)
$ws.Rows(39).Insert()

# Copy the "constant" columns (same for every data row in this sheet) from the
# row just below (which now holds what used to be row 39) so the new row
# matches the rest of the table: A Mercado ID, B Mercado, C Region,
# E Codreg, F Categoria ID, G Categoria, H Variedad, I Calidad,
# N Unidad de comercializacion, O Origen, Q Kg o Unidades, R Clasificacion.
$ws.Cells.Item(39, 1).Value  = $ws.Cells.Item(40, 1).Value2
$ws.Cells.Item(39, 2).Value  = $ws.Cells.Item(40, 2).Value2
$ws.Cells.Item(39, 3).Value  = $ws.Cells.Item(40, 3).Value2
$ws.Cells.Item(39, 5).Value  = $ws.Cells.Item(40, 5).Value2
$ws.Cells.Item(39, 6).Value  = $ws.Cells.Item(40, 6).Value2
$ws.Cells.Item(39, 7).Value  = $ws.Cells.Item(40, 7).Value2
$ws.Cells.Item(39, 8).Value  = $ws.Cells.Item(40, 8).Value2
$ws.Cells.Item(39, 9).Value  = $ws.Cells.Item(40, 9).Value2
$ws.Cells.Item(39, 14).Value = $ws.Cells.Item(40, 14).Value2
$ws.Cells.Item(39, 15).Value = $ws.Cells.Item(40, 15).Value2
$ws.Cells.Item(39, 17).Value = $ws.Cells.Item(40, 17).Value2
$ws.Cells.Item(39, 18).Value = $ws.Cells.Item(40, 18).Value2

# Row-specific values for the new record (D Fecha, J Volumen, K/L/M precios, P Precio $/Kg)
$ws.Cells.Item(39, 4).Value  = 44481
$ws.Cells.Item(39, 10).Value = 180
$ws.Cells.Item(39, 11).Value = 5500
$ws.Cells.Item(39, 12).Value = 5500
$ws.Cells.Item(39, 13).Value = 5500
$ws.Cells.Item(39, 16).Value = 153
